# "fix check for current holdings"
#
# Row 3 of the sheet (ZDB-ID 2238249-5, "The making of the modern world")
# had an incorrect current-holdings check: "Anzahl-FL-Bibliotheken" (count
# of libraries currently holding the title) was "0" and "FL-Bibliotheken"
# (the list of those libraries) was empty. Correct this to show the title
# is held by one library, with ID 547.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns X ("Anzahl-FL-Bibliotheken") and Y ("FL-Bibliotheken") store their
# values as text in this workbook (even when the text looks numeric), so
# force text formatting before assigning the values to avoid Excel silently
# re-typing them as numbers.
$range = $ws.Range("X3:Y3")
$range.NumberFormat = "@"

$ws.Range("X3").Value = "1"
$ws.Range("Y3").Value = "547"

# Restore the default cell style so no stray number-format style lingers on
# these cells (matches the rest of the sheet, which uses the default style).
$range.Style = "Normal"
